$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing rows 173-175 (B and D columns)
$ws.Range("B173").Value = 12213
$ws.Range("D173").Value = 1310

$ws.Range("B174").Value = 13721
$ws.Range("D174").Value = 1350

$ws.Range("B175").Value = 12526
$ws.Range("D175").Value = 1301

# Add new row 176 with a new month "01-07-2021".
# Build the text via a formula (string concatenation) in a scratch cell so
# Excel's "looks like a date" autoconversion never kicks in, then paste the
# resulting value (not the formula) into A176 and clean up the scratch cell.
$ws.Range("Z1").Formula = "=""01-0""&""7-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A176").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("B176").Value = 12994
$ws.Range("C176").Value = 2714
$ws.Range("D176").Value = 1375
$ws.Range("E176").Value = 287
